$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Verify Amazon User Login"
$ws.Range("C16").Value = "Passed"
